$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Format as text first so numeric-looking strings (e.g. "6.50", "37.30") keep their exact
# textual representation instead of being auto-converted to numbers by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "66.918.42"
$ws.Range("E2").Value = "  +2.03%  "
$ws.Range("D3").Value = "3.106.17"
$ws.Range("E3").Value = "  +5.09%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "581.19"
$ws.Range("E5").Value = "  +1.69%  "
$ws.Range("D6").Value = "172.97"
$ws.Range("E6").Value = "  +6.44%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "3.101.56"
$ws.Range("E8").Value = "  +5.05%  "
$ws.Range("E9").Value = "  +1.34%  "
$ws.Range("D10").Value = "6.50"
$ws.Range("E10").Value = "  -2.64%  "
$ws.Range("D11").Value = "0.157"
$ws.Range("E11").Value = "  +4.07%  "
$ws.Range("E12").Value = "  +4.98%  "
$ws.Range("E13").Value = "  +2.30%  "
$ws.Range("D14").Value = "37.30"
$ws.Range("E14").Value = "  +7.06%  "
$ws.Range("E15").Value = "  +0.16%  "
$ws.Range("D16").Value = "3.620.82"
$ws.Range("E16").Value = "  +5.06%  "
$ws.Range("D17").Value = "66.920.76"
$ws.Range("E18").Value = "  +2.22%  "
$ws.Range("D19").Value = "3.107.82"
$ws.Range("E19").Value = "  +5.20%  "
$ws.Range("D20").Value = "16.07"
$ws.Range("E20").Value = "  +0.76%  "
$ws.Range("D21").Value = "484.22"
$ws.Range("E21").Value = "  +8.60%  "
$ws.Range("E22").Value = "  +2.80%  "
$ws.Range("D23").Value = "7.53"
$ws.Range("E23").Value = "  +3.15%  "
$ws.Range("D24").Value = "84.14"
$ws.Range("E24").Value = "  +2.55%  "
$ws.Range("E25").Value = "  +5.95%  "
$ws.Range("D26").Value = "13.18"
$ws.Range("E26").Value = "  +7.29%  "
$ws.Range("D27").Value = "10.08"
$ws.Range("E27").Value = "  +0.15%  "
$ws.Range("D29").Value = "7.99"
$ws.Range("E29").Value = "  -1.49%  "
$ws.Range("D30").Value = "2.39"
$ws.Range("E30").Value = "  -5.41%  "
$ws.Range("D31").Value = "2.70"
$ws.Range("E31").Value = "  +3.78%  "
$ws.Range("E32").Value = "  -0.30%  "
$ws.Range("D33").Value = "28.89"
$ws.Range("E33").Value = "  +6.33%  "
$ws.Range("D34").Value = "0.115"
$ws.Range("E34").Value = "  +1.05%  "
$ws.Range("E35").Value = "  +0.09%  "
$ws.Range("E36").Value = "  +3.45%  "
$ws.Range("D37").Value = "0.998"
$ws.Range("E37").Value = "  +2.33%  "
$ws.Range("D38").Value = "48.07"
$ws.Range("E38").Value = "  +4.10%  "
$ws.Range("E39").Value = "  +6.94%  "
$ws.Range("E40").Value = "  +4.79%  "
$ws.Range("D41").Value = "50.13"
$ws.Range("E41").Value = "  +2.12%  "
$ws.Range("E42").Value = "  -0.07%  "
$ws.Range("D43").Value = "8.67"
$ws.Range("E43").Value = "  +1.62%  "
$ws.Range("D44").Value = "2.82"
$ws.Range("E44").Value = "  -0.08%  "
$ws.Range("D45").Value = "0.0363"
$ws.Range("E45").Value = "  +3.21%  "
$ws.Range("D46").Value = "2.836.65"
$ws.Range("E46").Value = "  +5.81%  "
$ws.Range("D47").Value = "384.14"
$ws.Range("E47").Value = "  -0.18%  "
$ws.Range("D48").Value = "135.11"
$ws.Range("E48").Value = "  +1.68%  "
$ws.Range("E49").Value = "  +0.00%  "
$ws.Range("D50").Value = "24.92"
$ws.Range("E50").Value = "  +4.26%  "
$ws.Range("E51").Value = "  +3.09%  "

# Restore the original (default) cell style now that the text values are locked in,
# so the saved file does not carry a stray Text number-format style.
$ws.Range("D2:E51").Style = "Normal"
